$d = $word.ActiveDocument

# Each of these merge fields is currently written as three runs:
#   run1 = "${"   run2 = <name> (wrapped in w:proofErr spellStart/spellEnd)   run3 = "}"
# The edit collapses each occurrence down to a single run whose text is
# "${<newName>}" and drops the w:proofErr markers. For "departmentFull" the
# field name itself also changes (to "support"); for the others the visible
# text ends up the same as before, only the run/proofErr structure changes.
$replacements = @(
    @{ Old = "departmentFull"; New = "support" },
    @{ Old = "ilos";           New = "ilos" },
    @{ Old = "budgetSource";   New = "budgetSource" },
    @{ Old = "sig_cscp";       New = "sig_cscp" },
    @{ Old = "sig_csca";       New = "sig_csca" },
    @{ Old = "sig_sscp";       New = "sig_sscp" },
    @{ Old = "sig_dean";       New = "sig_dean" }
)

foreach ($item in $replacements) {
    $oldName = $item.Old
    $newName = $item.New

    $rng = $d.Content
    $found = $rng.Find.Execute($oldName, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $oldName"
        continue
    }

    $start = $rng.Start
    $end = $rng.End

    # Widen the match by one character on each side so the range swallows
    # the neighbouring "${" and "}" runs too, crossing the w:proofErr
    # spellStart/spellEnd boundaries that surround the middle run.
    $wideStart = $start - 2
    $wideEnd = $end + 1

    $finalText = "`${" + $newName + "}"

    $wide = $d.Range($wideStart, $wideEnd)

    # Word's Range.Text setter is a no-op when the replacement text is
    # byte-for-byte identical to what's already there (which happens for
    # every field here except departmentFull/support, since the visible
    # text "${name}" is not actually changing). Force the run-merge (and
    # the associated w:proofErr cleanup) to happen by first swapping in a
    # guaranteed-different placeholder, then writing the real final text.
    $placeholder = "___TMP_PLACEHOLDER___"
    $wide.Text = $placeholder

    $wide2 = $d.Range($wideStart, $wideStart + $placeholder.Length)
    $wide2.Text = $finalText
}
